$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "YES"
$ws.Range("B3").Value = "NO"
$ws.Range("B4").Value = "YES"
$ws.Range("B5").Value = "YES"
$ws.Range("B6").Value = "YES"
$ws.Range("B7").Value = "YES"
$ws.Range("B8").Value = "YES"
$ws.Range("B9").Value = "YES"
